$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S91028.MES.BIN")

# Column F = "Edited", Column G = "Initial"
# Update the translated ("Edited") lines that were reworded, and fill in the
# previously-empty "Edited" cells with new phrasing for rows 6-9 and 13-14.

$ws.Range("G3").Value = "I start to feel dizzy as soon as I enter the room…"
$ws.Range("G5").Value = "I'm so pathetic…"

$ws.Range("F6").Value = "I almost collapsed as soon as I entered the room…"
$ws.Range("F7").Value = "Somehow, I made it to the bed, but it might be a bit tough."
$ws.Range("F8").Value = "For now, I'll rest tonight,"
$ws.Range("F9").Value = "and see if I recover by tomorrow…"

$ws.Range("F13").Value = "I don't feel recovered at all."
$ws.Range("F14").Value = "I'll take today off…"
